# Add a new "SlideCollection.add" example row to the Snippets table,
# showing how to use the SlideCollection.add method (powerpoint-add-slides).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data lives in an Excel Table ("Snippets") - grow it by one row so the
# table range / autofilter / sheet dimension all expand together.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# New row lands at A3:E3 (columns: Class, Method/Prop/Rel Name,
# Member ID (methods only), SnippetIdIntheYAMLFile, MethodNameInTheSnippet).
$ws.Range("A3").Value = "SlideCollection"
$ws.Range("B3").Value = "add"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "powerpoint-add-slides"
$ws.Range("E3").Value = "addSlide"

# Match the explicit "General" number format applied to the text columns
# (the numeric "Member ID" column C is left with the default/no style).
$ws.Range("A3:B3").NumberFormat = "General"
$ws.Range("D3:E3").NumberFormat = "General"

# Column C (Member ID) got a bit wider to fit the new content.
$ws.Columns.Item(3).ColumnWidth = 29.498697916666668

# Move the active selection down past the new row, like after typing it in.
$null = $ws.Range("E5").Select()
